$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 970.4375
$ws.Range("I2").Value = 428.8
$ws.Range("J2").Value = 1873.1666
$ws.Range("K2").Value = 428.8
$ws.Range("L2").Value = 1873.1666
$ws.Range("M2").Value = -315.8
$ws.Range("N2").Value = -2099.1666
$ws.Range("H12").Value = 888.8333
$ws.Range("I12").Value = 867.2
$ws.Range("K12").Value = 867.2
$ws.Range("M12").Value = -697.2
$ws.Range("H98").Value = 1985
$ws.Range("I98").Value = 1166.1818
$ws.Range("K98").Value = 1166.1818
$ws.Range("M98").Value = 331.8181999999999
$ws.Range("H122").Value = 1985
$ws.Range("I122").Value = 1166.1818
$ws.Range("K122").Value = 3498.5454
$ws.Range("M122").Value = -1048.5454
$ws.Range("H129").Value = 876.7
$ws.Range("I129").Value = 876.7
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2630.1
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2369.9
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 18321.904
$ws.Range("I137").Value = 23746.732
$ws.Range("K137").Value = 71240.196
$ws.Range("M137").Value = -68690.196
$ws.Range("H138").Value = 3899.1162
$ws.Range("J138").Value = 5600.5884
$ws.Range("L138").Value = 16801.7652
$ws.Range("N138").Value = -27081.7652
$ws.Range("H141").Value = 1479.85
$ws.Range("I141").Value = 1479.85
$ws.Range("K141").Value = 4439.549999999999
$ws.Range("M141").Value = 740.4500000000007

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2043.125
$ws.Range("I2").Value = 1766.44
$ws.Range("K2").Value = 1766.44
$ws.Range("M2").Value = -1653.44
$ws.Range("H4").Value = 637.25
$ws.Range("I4").Value = 574.6667
$ws.Range("J4").Value = 825
$ws.Range("K4").Value = 574.6667
$ws.Range("L4").Value = 825
$ws.Range("M4").Value = -458.6667
$ws.Range("N4").Value = -1057
$ws.Range("H18").Value = 20000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H22").Value = 5671.25
$ws.Range("J22").Value = 6728.6665
$ws.Range("L22").Value = 6728.6665
$ws.Range("N22").Value = -7326.6665
$ws.Range("H45").Value = 3346.2632
$ws.Range("I45").Value = 1944.5927
$ws.Range("K45").Value = 1944.5927
$ws.Range("M45").Value = -1567.5927
$ws.Range("H74").Value = 241291.52
$ws.Range("I74").Value = 261664.73
$ws.Range("J74").Value = 6999.5
$ws.Range("K74").Value = 261664.73
$ws.Range("L74").Value = 6999.5
$ws.Range("M74").Value = -260790.73
$ws.Range("N74").Value = -8747.5
$ws.Range("H77").Value = 241291.52
$ws.Range("I77").Value = 261664.73
$ws.Range("J77").Value = 6999.5
$ws.Range("K77").Value = 1308323.65
$ws.Range("L77").Value = 34997.5
$ws.Range("M77").Value = -1303955.65
$ws.Range("N77").Value = -43733.5
$ws.Range("H116").Value = 2043.125
$ws.Range("I116").Value = 1766.44
$ws.Range("K116").Value = 1766.44
$ws.Range("M116").Value = 527.5599999999999
$ws.Range("H122").Value = 1985.7222
$ws.Range("I122").Value = 1744.258
$ws.Range("J122").Value = 3482.8
$ws.Range("K122").Value = 5232.774
$ws.Range("L122").Value = 10448.4
$ws.Range("M122").Value = -2782.774
$ws.Range("N122").Value = -15348.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2043.125
$ws.Range("I3").Value = 1766.44
$ws.Range("K3").Value = 1766.44
$ws.Range("M3").Value = -1652.44
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H107").Value = 20042.268
$ws.Range("I107").Value = 25649.137
$ws.Range("K107").Value = 25649.137
$ws.Range("M107").Value = -23729.137

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 67971
$ws.Range("J116").Value = 67971
$ws.Range("L116").Value = 67971
$ws.Range("N116").Value = -77149
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2700
$ws.Range("I48").Value = 900
$ws.Range("K48").Value = 2700
$ws.Range("M48").Value = -2450
$ws.Range("H121").Value = 76718.64
$ws.Range("I121").Value = 355174.66
$ws.Range("J121").Value = 776.0909
$ws.Range("K121").Value = 1065523.98
$ws.Range("L121").Value = 2328.2727
$ws.Range("M121").Value = -1064213.98
$ws.Range("N121").Value = -4948.2727
$ws.Range("H129").Value = 4202.095
$ws.Range("J129").Value = 4652
$ws.Range("L129").Value = 13956
$ws.Range("N129").Value = -23956

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 26127
$ws.Range("J20").Value = 26127
$ws.Range("L20").Value = 26127
$ws.Range("N20").Value = -26617
$ws.Range("H70").Value = 19495.592
$ws.Range("J70").Value = 33847.5
$ws.Range("L70").Value = 33847.5
$ws.Range("N70").Value = -34387.5
$ws.Range("H73").Value = 19495.592
$ws.Range("J73").Value = 33847.5
$ws.Range("L73").Value = 33847.5
$ws.Range("N73").Value = -35719.5
$ws.Range("H80").Value = 15399.8
$ws.Range("I80").Value = 6749.75
$ws.Range("K80").Value = 6749.75
$ws.Range("M80").Value = -5751.75
$ws.Range("H83").Value = 15399.8
$ws.Range("I83").Value = 6749.75
$ws.Range("K83").Value = 33748.75
$ws.Range("M83").Value = -28756.75
$ws.Range("H102").Value = 28015.63
$ws.Range("I102").Value = 32857
$ws.Range("K102").Value = 32857
$ws.Range("M102").Value = -31235
$ws.Range("H119").Value = 59999
$ws.Range("J119").Value = 59999
$ws.Range("L119").Value = 59999
$ws.Range("N119").Value = -69675

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3652.6428
$ws.Range("I7").Value = 3472.077
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 3472.077
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -3360.077
$ws.Range("N7").Value = -6224
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H40").Value = 3079.8
$ws.Range("I40").Value = 3155.3333
$ws.Range("K40").Value = 3155.3333
$ws.Range("M40").Value = -3019.3333
$ws.Range("H122").Value = 16999.23
$ws.Range("I122").Value = 16999.23
$ws.Range("K122").Value = 50997.69
$ws.Range("M122").Value = -48547.69
$ws.Range("H126").Value = 3652.6428
$ws.Range("I126").Value = 3472.077
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 10416.231
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -7946.231
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 2481.2307
$ws.Range("I132").Value = 2380.48
$ws.Range("K132").Value = 7141.440000000001
$ws.Range("M132").Value = -4611.440000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1864.15
$ws.Range("I113").Value = 1181.6364
$ws.Range("K113").Value = 3544.9092
$ws.Range("M113").Value = -1374.9092
$ws.Range("H122").Value = 14790951
$ws.Range("I122").Value = 21863004
$ws.Range("K122").Value = 65589012
$ws.Range("M122").Value = -65586562
$ws.Range("H126").Value = 180809.36
$ws.Range("I126").Value = 1849.9546
$ws.Range("J126").Value = 836993.8
$ws.Range("K126").Value = 5549.8638
$ws.Range("L126").Value = 2510981.4
$ws.Range("M126").Value = -3079.8638
$ws.Range("N126").Value = -2515921.4
